$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.03863476054525
$ws.Range("D2").Value = 0.008436957273197265
$ws.Range("E2").Value = 0.5149716421402282
$ws.Range("F2").Value = 0.3786870749330831
$ws.Range("G2").Value = 0.2471652642554574
$ws.Range("H2").Value = 0.3599890557817389
$ws.Range("L2").Value = 0.2346600967121191
$ws.Range("N2").Value = 1.536221915440194
$ws.Range("O2").Value = 1.129660440777883
$ws.Range("B3").Value = 0.9896902276415176
$ws.Range("D3").Value = 0.007347719088468097
$ws.Range("E3").Value = 0.5152697082145643
$ws.Range("F3").Value = 0.3635543690312915
$ws.Range("G3").Value = 0.2326764080411863
$ws.Range("H3").Value = 0.3566495377782672
$ws.Range("L3").Value = 0.2075716459574011
$ws.Range("N3").Value = 1.517045385650164
$ws.Range("O3").Value = 1.091669330473252
$ws.Range("B4").Value = 0.9599708089972694
$ws.Range("D4").Value = 0.006677667548334654
$ws.Range("E4").Value = 0.5156377326705623
$ws.Range("F4").Value = 0.3545305629466426
$ws.Range("G4").Value = 0.2239665933931008
$ws.Range("H4").Value = 0.35481175322154
$ws.Range("L4").Value = 0.1909186671224319
$ws.Range("N4").Value = 1.505885600044309
$ws.Range("O4").Value = 1.069172406452282
$ws.Range("B5").Value = 0.9479445523829213
$ws.Range("D5").Value = 0.006404317409852922
$ws.Range("E5").Value = 0.5158345712925154
$ws.Range("F5").Value = 0.3509205666361481
$ws.Range("G5").Value = 0.2204640078322342
$ws.Range("H5").Value = 0.354116347192786
$ws.Range("L5").Value = 0.1841277135008568
$ws.Range("N5").Value = 1.501493579141538
$ws.Range("O5").Value = 1.060213128225797
$ws.Range("B6").Value = 0.9459527430903165
$ws.Range("D6").Value = 0.006358910281953456
$ws.Range("E6").Value = 0.5158700968585848
$ws.Range("F6").Value = 0.3503251915554841
$ws.Range("G6").Value = 0.2198852247870349
$ws.Range("H6").Value = 0.354004107861229
$ws.Range("L6").Value = 0.1829998081578879
$ws.Range("N6").Value = 1.500773723218629
$ws.Range("O6").Value = 1.058738026431996
$ws.Range("B7").Value = 0.9598082746578882
$ws.Range("D7").Value = 0.006673982239131959
$ws.Range("E7").Value = 0.5156401970915248
$ws.Range("F7").Value = 0.3544816049311734
$ws.Range("G7").Value = 0.2239191672558576
$ws.Range("H7").Value = 0.3548021580595844
$ws.Range("L7").Value = 0.1908271005115836
$ws.Range("N7").Value = 1.505825735876215
$ws.Range("O7").Value = 1.06905073502304
$ws.Range("B8").Value = 1.021690200742086
$ws.Range("D8").Value = 0.008061658204475464
$ws.Range("E8").Value = 0.5150362210099715
$ws.Range("F8").Value = 0.3734137158926387
$ws.Range("G8").Value = 0.242130731686288
$ws.Range("H8").Value = 0.3587934571548033
$ws.Range("L8").Value = 0.225324522076221
$ws.Range("N8").Value = 1.529482927875094
$ws.Range("O8").Value = 1.116388721575078
$ws.Range("B9").Value = 1.145644080625999
$ws.Range("D9").Value = 0.01077231353330177
$ws.Range("E9").Value = 0.5153050863985058
$ws.Range("F9").Value = 0.4126684142082411
$ws.Range("G9").Value = 0.2793308307415003
$ws.Range("H9").Value = 0.3683079259542836
$ws.Range("L9").Value = 0.2927946133380317
$ws.Range("N9").Value = 1.580706909472894
$ws.Range("O9").Value = 1.215819927545766
$ws.Range("B10").Value = 1.238261259267802
$ws.Range("D10").Value = 0.01275674202575772
$ws.Range("E10").Value = 0.5163693156351812
$ws.Range("F10").Value = 0.442816398869553
$ws.Range("G10").Value = 0.3075828245487742
$ws.Range("H10").Value = 0.3763282297441464
$ws.Range("L10").Value = 0.3422387495918713
$ws.Range("N10").Value = 1.621232702405536
$ws.Range("O10").Value = 1.292930142125414
$ws.Range("B11").Value = 1.280723766144604
$ws.Range("D11").Value = 0.01365784338958065
$ws.Range("E11").Value = 0.5170376187554808
$ws.Range("F11").Value = 0.4568176881553114
$ws.Range("G11").Value = 0.3206386259047918
$ws.Range("H11").Value = 0.3802008831299162
$ws.Range("L11").Value = 0.3647014958215493
$ws.Range("N11").Value = 1.640285799000338
$ws.Range("O11").Value = 1.328898443853916
$ws.Range("B12").Value = 1.29684987764432
$ws.Range("D12").Value = 0.01399881819674675
$ws.Range("E12").Value = 0.5173168306369504
$ws.Range("F12").Value = 0.4621609571321841
$ws.Range("G12").Value = 0.3256120169456977
$ws.Range("H12").Value = 0.3816995921180393
$ws.Range("L12").Value = 0.3732029073025558
$ws.Range("N12").Value = 1.64758852066646
$ws.Range("O12").Value = 1.342647160183475
$ws.Range("B13").Value = 1.293374780284637
$ws.Range("D13").Value = 0.01392539471449794
$ws.Range("E13").Value = 0.5172555402554693
$ws.Range("F13").Value = 0.4610083517972328
$ws.Range("G13").Value = 0.3245395966613103
$ws.Range("H13").Value = 0.3813753856435653
$ws.Range("L13").Value = 0.3713721948245166
$ws.Range("N13").Value = 1.646011864571278
$ws.Range("O13").Value = 1.339680417967713
$ws.Range("B14").Value = 1.282049545442305
$ws.Range("D14").Value = 0.01368590076567244
$ws.Range("E14").Value = 0.5170600677707782
$ws.Range("F14").Value = 0.4572564544836837
$ws.Range("G14").Value = 0.3210471989604144
$ws.Range("H14").Value = 0.3803235371995299
$ws.Range("L14").Value = 0.3654010098690321
$ws.Range("N14").Value = 1.640884846372671
$ws.Range("O14").Value = 1.330026986260037
$ws.Range("B15").Value = 1.275118534489195
$ws.Range("D15").Value = 0.01353917032933794
$ws.Range("E15").Value = 0.5169437292183545
$ws.Range("F15").Value = 0.4549636866689184
$ws.Range("G15").Value = 0.318911842681004
$ws.Range("H15").Value = 0.3796834450276094
$ws.Range("L15").Value = 0.3617428547322561
$ws.Range("N15").Value = 1.637755791969965
$ws.Range("O15").Value = 1.324130696536713
$ws.Range("B16").Value = 1.235492795718926
$ws.Range("D16").Value = 0.01269781860538188
$ws.Range("E16").Value = 0.5163293171134029
$ws.Range("F16").Value = 0.4419071567337625
$ws.Range("G16").Value = 0.3067337101222876
$ws.Range("H16").Value = 0.3760796522959424
$ws.Range("L16").Value = 0.3407701187480541
$ws.Range("N16").Value = 1.61999987350481
$ws.Range("O16").Value = 1.290597469358318
$ws.Range("B17").Value = 1.211267592253932
$ws.Range("D17").Value = 0.01218124682207389
$ws.Range("E17").Value = 0.5159993222552899
$ws.Range("F17").Value = 0.4339708894050887
$ws.Range("G17").Value = 0.2993151203060762
$ws.Range("H17").Value = 0.3739262472385292
$ws.Range("L17").Value = 0.3278960907086059
$ws.Range("N17").Value = 1.609264618155748
$ws.Range("O17").Value = 1.270254174459325
$ws.Range("B18").Value = 1.197365065181998
$ws.Range("D18").Value = 0.01188397657399065
$ws.Range("E18").Value = 0.5158268743520509
$ws.Range("F18").Value = 0.4294331529382873
$ws.Range("G18").Value = 0.2950673211381201
$ws.Range("H18").Value = 0.3727087672600646
$ws.Range("L18").Value = 0.3204885291794142
$ws.Range("N18").Value = 1.603148190019425
$ws.Range("O18").Value = 1.258637050714782
$ws.Range("B19").Value = 1.192663293525527
$ws.Range("D19").Value = 0.01178330047947185
$ws.Range("E19").Value = 0.5157714779664886
$ws.Range("F19").Value = 0.4279013884203948
$ws.Range("G19").Value = 0.2936323790067945
$ws.Range("H19").Value = 0.3723001742043408
$ws.Range("L19").Value = 0.3179799965706991
$ws.Range("N19").Value = 1.601087303263512
$ws.Range("O19").Value = 1.254718077434688
$ws.Range("B20").Value = 1.213843188173087
$ws.Range("D20").Value = 0.0122362526209443
$ws.Range("E20").Value = 0.5160326569508129
$ws.Range("F20").Value = 0.4348129244716006
$ws.Range("G20").Value = 0.3001028560937584
$ws.Range("H20").Value = 0.3741532970626196
$ws.Range("L20").Value = 0.3292668422222391
$ws.Range("N20").Value = 1.610401388217838
$ws.Range("O20").Value = 1.272411077954501
$ws.Range("B21").Value = 1.285374788238414
$ws.Range("D21").Value = 0.01375625295535343
$ws.Range("E21").Value = 0.5171167760872137
$ws.Range("F21").Value = 0.4583573569118755
$ws.Range("G21").Value = 0.3220722007222889
$ws.Range("H21").Value = 0.3806316163049814
$ws.Range("L21").Value = 0.3671550241861326
$ws.Range("N21").Value = 1.642388403702682
$ws.Range("O21").Value = 1.332858949130269
$ws.Range("B22").Value = 1.332395272710528
$ws.Range("D22").Value = 0.01474817910423809
$ws.Range("E22").Value = 0.5179775154936621
$ws.Range("F22").Value = 0.4739856596437448
$ws.Range("G22").Value = 0.3366021099320875
$ws.Range("H22").Value = 0.3850533622054826
$ws.Range("L22").Value = 0.3918892729243737
$ws.Range("N22").Value = 1.663804607876614
$ws.Range("O22").Value = 1.373113076142573
$ws.Range("B23").Value = 1.307275152946033
$ws.Range("D23").Value = 0.01421891135288433
$ws.Range("E23").Value = 0.5175043094843161
$ws.Range("F23").Value = 0.46562250775591
$ws.Range("G23").Value = 0.328831472855498
$ws.Range("H23").Value = 0.3826762162462671
$ws.Range("L23").Value = 0.3786908474932034
$ws.Range("N23").Value = 1.652327995009927
$ws.Range("O23").Value = 1.351560172271263
$ws.Range("B24").Value = 1.212678683723368
$ws.Range("D24").Value = 0.01221138538932109
$ws.Range("E24").Value = 0.5160175325133949
$ws.Range("F24").Value = 0.4344321627684451
$ws.Range("G24").Value = 0.2997466670321529
$ws.Range("H24").Value = 0.3740505838496375
$ws.Range("L24").Value = 0.3286471444639574
$ws.Range("N24").Value = 1.609887281781766
$ws.Range("O24").Value = 1.271435697357703
$ws.Range("B25").Value = 1.111836574486034
$ws.Range("D25").Value = 0.01004020773459757
$ws.Range("E25").Value = 0.5150787647469741
$ws.Range("F25").Value = 0.4018200518247284
$ws.Range("G25").Value = 0.2691064982479503
$ws.Range("H25").Value = 0.3655532158409329
$ws.Range("L25").Value = 0.2745631029248727
$ws.Range("N25").Value = 1.56633779501081
$ws.Range("O25").Value = 1.18821091226701
